$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# --- Update F-column timestamps on the "data" sheet (rows 2..99) ---
$timestamps = @(
  "2021-10-05 14:33:11.987927",
  "2021-10-05 14:33:11.987936",
  "2021-10-05 14:33:11.987939",
  "2021-10-05 14:33:11.987942",
  "2021-10-05 14:33:11.987945",
  "2021-10-05 14:33:11.987948",
  "2021-10-05 14:33:11.987951",
  "2021-10-05 14:33:11.987954",
  "2021-10-05 14:33:11.987956",
  "2021-10-05 14:33:11.987959",
  "2021-10-05 14:33:11.987962",
  "2021-10-05 14:33:11.987964",
  "2021-10-05 14:33:11.987967",
  "2021-10-05 14:33:11.987970",
  "2021-10-05 14:33:11.987972",
  "2021-10-05 14:33:11.987975",
  "2021-10-05 14:33:11.987978",
  "2021-10-05 14:33:11.987982",
  "2021-10-05 14:33:11.987985",
  "2021-10-05 14:33:11.987987",
  "2021-10-05 14:33:11.987990",
  "2021-10-05 14:33:11.987992",
  "2021-10-05 14:33:11.987995",
  "2021-10-05 14:33:11.987998",
  "2021-10-05 14:33:11.988001",
  "2021-10-05 14:33:11.988004",
  "2021-10-05 14:33:11.988006",
  "2021-10-05 14:33:11.988009",
  "2021-10-05 14:33:11.988012",
  "2021-10-05 14:33:11.988014",
  "2021-10-05 14:33:11.988017",
  "2021-10-05 14:33:11.988020",
  "2021-10-05 14:33:11.988023",
  "2021-10-05 14:33:11.988026",
  "2021-10-05 14:33:11.988028",
  "2021-10-05 14:33:11.988031",
  "2021-10-05 14:33:11.988034",
  "2021-10-05 14:33:11.988036",
  "2021-10-05 14:33:11.988039",
  "2021-10-05 14:33:11.988042",
  "2021-10-05 14:33:11.988045",
  "2021-10-05 14:33:11.988047",
  "2021-10-05 14:33:11.988050",
  "2021-10-05 14:33:11.988053",
  "2021-10-05 14:33:11.988055",
  "2021-10-05 14:33:11.988058",
  "2021-10-05 14:33:11.988061",
  "2021-10-05 14:33:11.988063",
  "2021-10-05 14:33:11.988066",
  "2021-10-05 14:33:11.988069",
  "2021-10-05 14:33:11.988071",
  "2021-10-05 14:33:11.988074",
  "2021-10-05 14:33:11.988077",
  "2021-10-05 14:33:11.988079",
  "2021-10-05 14:33:11.988082",
  "2021-10-05 14:33:11.988084",
  "2021-10-05 14:33:11.988087",
  "2021-10-05 14:33:11.988090",
  "2021-10-05 14:33:11.988092",
  "2021-10-05 14:33:11.988095",
  "2021-10-05 14:33:11.988097",
  "2021-10-05 14:33:11.988100",
  "2021-10-05 14:33:11.988103",
  "2021-10-05 14:33:11.988105",
  "2021-10-05 14:33:11.988109",
  "2021-10-05 14:33:11.988112",
  "2021-10-05 14:33:11.988115",
  "2021-10-05 14:33:11.988117",
  "2021-10-05 14:33:11.988120",
  "2021-10-05 14:33:11.988123",
  "2021-10-05 14:33:11.988125",
  "2021-10-05 14:33:11.988128",
  "2021-10-05 14:33:11.988131",
  "2021-10-05 14:33:11.988133",
  "2021-10-05 14:33:11.988136",
  "2021-10-05 14:33:11.988139",
  "2021-10-05 14:33:11.988144",
  "2021-10-05 14:33:11.988147",
  "2021-10-05 14:33:11.988150",
  "2021-10-05 14:33:11.988152",
  "2021-10-05 14:33:11.988155",
  "2021-10-05 14:33:11.988158",
  "2021-10-05 14:33:11.988160",
  "2021-10-05 14:33:11.988163",
  "2021-10-05 14:33:11.988165",
  "2021-10-05 14:33:11.988168",
  "2021-10-05 14:33:11.988171",
  "2021-10-05 14:33:11.988174",
  "2021-10-05 14:33:11.988176",
  "2021-10-05 14:33:11.988179",
  "2021-10-05 14:33:11.988182",
  "2021-10-05 14:33:11.988184",
  "2021-10-05 14:33:11.988188",
  "2021-10-05 14:33:11.988191",
  "2021-10-05 14:33:11.988194",
  "2021-10-05 14:33:11.988196",
  "2021-10-05 14:33:11.988199",
  "2021-10-05 14:33:11.988202"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $data.Cells.Item($row, 6).Value = $timestamps[$i]
}

# --- Add the "metadata" worksheet after "data" ---
$ws = $wb.Worksheets.Add($null, $data)
$ws.Name = "metadata"

# Copy the bold/border/center header style used on "data"!A2 (and the
# header row) onto the new header row + the A2 index cell.
$data.Range("A2").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Ataxia - adult onset"
$ws.Range("C2").Value = 268

# data_version must stay a text string ("0.143"), not become a float.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0.143"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "2021-10-04T04:59:19.151898Z"
$ws.Range("F2").Value = "2021-10-05 14:33:11.984452"
$ws.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/268/?format=json"

$ws.Activate()
$ws.Range("A1").Select()

Write-Output "edit complete"
